# "Fruta / hortaliza, semanal" - weekly refresh of the Achicoria price
# series for Vega Modelo de Temuco.
#
# The new weekly observation is inserted as row 18 (pushing every
# existing record from row 18 down to row 19, all the way through the
# former last row 143, which becomes row 144). The newly inserted row
# holds this week's data point; every other row keeps its original
# content, just shifted down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; this shifts rows 18:143 down to
# 19:144 and keeps all of their existing values/formatting intact.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with this week's record.
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 45163
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 100112010
$ws.Range("G18").Value = "Achicoria"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 10000
$ws.Range("N18").Value = "$/caja 18 unidades"
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 556
$ws.Range("Q18").Value = 18
$ws.Range("R18").Value = "Hortaliza"
